$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Sheet, $Row, $Col, $Val) {
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

Set-TextCell $ws 2 4 '26.651.49'
Set-TextCell $ws 2 5 '  +4.30%  '
Set-TextCell $ws 3 4 '1.752.97'
Set-TextCell $ws 3 5 '  +4.96%  '
Set-TextCell $ws 4 4 '0.9962'
Set-TextCell $ws 5 4 '247.99'
Set-TextCell $ws 5 5 '  +3.83%  '
Set-TextCell $ws 6 4 '0.9968'
Set-TextCell $ws 6 5 '  -0.30%  '
Set-TextCell $ws 7 4 '0.4831'
Set-TextCell $ws 7 5 '  +0.39%  '
Set-TextCell $ws 8 4 '0.2715'
Set-TextCell $ws 8 5 '  +3.21%  '
Set-TextCell $ws 9 4 '0.06254'
Set-TextCell $ws 9 5 '  +1.03%  '
Set-TextCell $ws 10 4 '1.740.48'
Set-TextCell $ws 10 5 '  +4.19%  '
Set-TextCell $ws 11 4 '0.07123'
Set-TextCell $ws 11 5 '  +1.61%  '
Set-TextCell $ws 12 4 '16.02'
Set-TextCell $ws 12 5 '  +7.49%  '
Set-TextCell $ws 13 4 '0.6258'
Set-TextCell $ws 13 5 '  +5.85%  '
Set-TextCell $ws 14 4 '4.513'
Set-TextCell $ws 14 5 '  +2.90%  '
Set-TextCell $ws 15 4 '77.39'
Set-TextCell $ws 15 5 '  +2.90%  '
Set-TextCell $ws 16 4 '0.9968'
Set-TextCell $ws 16 5 '  -0.29%  '
Set-TextCell $ws 17 4 '26.653.92'
Set-TextCell $ws 17 5 '  +4.36%  '
Set-TextCell $ws 18 4 '0.9972'
Set-TextCell $ws 18 5 '  -0.20%  '
Set-TextCell $ws 19 4 '0.000006916'
Set-TextCell $ws 19 5 '  +2.04%  '
Set-TextCell $ws 20 4 '11.76'
Set-TextCell $ws 20 5 '  +2.75%  '
Set-TextCell $ws 21 4 '1.962.80'
Set-TextCell $ws 21 5 '  +4.33%  '
Set-TextCell $ws 22 4 '4.637'
Set-TextCell $ws 22 5 '  +4.23%  '
Set-TextCell $ws 23 4 '8.882'
Set-TextCell $ws 23 5 '  +1.57%  '
Set-TextCell $ws 24 4 '5.372'
Set-TextCell $ws 24 5 '  +1.60%  '
Set-TextCell $ws 25 4 '136.36'
Set-TextCell $ws 25 5 '  -0.27%  '
Set-TextCell $ws 26 4 '15.48'
Set-TextCell $ws 26 5 '  +2.85%  '
Set-TextCell $ws 27 4 '1.843'
Set-TextCell $ws 27 5 '  +6.73%  '
Set-TextCell $ws 28 4 '1.412'
Set-TextCell $ws 28 5 '  +1.51%  '
Set-TextCell $ws 29 4 '108.03'
Set-TextCell $ws 29 5 '  +3.02%  '
Set-TextCell $ws 30 4 '4.032'
Set-TextCell $ws 30 5 '  +1.45%  '
Set-TextCell $ws 31 4 '3.778'
Set-TextCell $ws 31 5 '  +3.32%  '
Set-TextCell $ws 32 4 '0.07911'
Set-TextCell $ws 32 5 '  +1.31%  '
Set-TextCell $ws 33 4 '0.04592'
Set-TextCell $ws 33 5 '  +8.11%  '
Set-TextCell $ws 34 4 '2.607'
Set-TextCell $ws 34 5 '  -0.09%  '
Set-TextCell $ws 35 4 '1.007'
Set-TextCell $ws 35 5 '  +5.70%  '
Set-TextCell $ws 36 4 '0.6367'
Set-TextCell $ws 36 5 '  +4.38%  '
Set-TextCell $ws 37 4 '0.9544'
Set-TextCell $ws 37 5 '  +10.99%  '
Set-TextCell $ws 38 4 '114.89'
Set-TextCell $ws 38 5 '  +19.72%  '
Set-TextCell $ws 39 4 '2.500'
Set-TextCell $ws 39 5 '  -3.52%  '
Set-TextCell $ws 40 4 '2.003'
Set-TextCell $ws 40 5 '  +7.74%  '
Set-TextCell $ws 41 5 '  +0.17%  '
Set-TextCell $ws 42 4 '5.752'
Set-TextCell $ws 42 5 '  +18.75%  '
Set-TextCell $ws 43 4 '0.01514'
Set-TextCell $ws 43 5 '  +2.11%  '
Set-TextCell $ws 44 4 '0.3932'
Set-TextCell $ws 44 5 '  +4.06%  '
Set-TextCell $ws 45 4 '6.798'
Set-TextCell $ws 45 5 '  +9.40%  '
Set-TextCell $ws 46 4 '0.1209'
Set-TextCell $ws 46 5 '  +7.91%  '
Set-TextCell $ws 47 4 '0.05332'
Set-TextCell $ws 47 5 '  +1.58%  '
Set-TextCell $ws 48 4 '7.967'
Set-TextCell $ws 48 5 '  +8.31%  '
Set-TextCell $ws 49 4 '30.98'
Set-TextCell $ws 49 5 '  +3.85%  '
Set-TextCell $ws 50 4 '0.3471'
Set-TextCell $ws 50 5 '  +3.90%  '
Set-TextCell $ws 51 2 'Aave'
Set-TextCell $ws 51 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 51 4 '51.97'
Set-TextCell $ws 51 5 '  +3.81%  '
